$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.615.94"
$ws.Range("E2").Value = "  +0.59%  "

# Row 3
$ws.Range("D3").Value = "2.517.73"
$ws.Range("E3").Value = "  +0.37%  "

# Row 4
$ws.Range("E4").Value = "  +0.11%  "

# Row 5
$ws.Range("D5").Value = "314.98"
$ws.Range("E5").Value = "  +2.64%  "

# Row 6
$ws.Range("D6").Value = "95.25"
$ws.Range("E6").Value = "  -1.05%  "

# Row 7
$ws.Range("D7").Value = "0.574"
$ws.Range("E7").Value = "  -1.92%  "

# Row 8
$ws.Range("E8").Value = "  -0.06%  "

# Row 9
$ws.Range("D9").Value = "0.532"
$ws.Range("E9").Value = "  -0.72%  "

# Row 10
$ws.Range("D10").Value = "35.75"
$ws.Range("E10").Value = "  -1.55%  "

# Row 11
$ws.Range("D11").Value = "0.0809"
$ws.Range("E11").Value = "  -0.34%  "

# Row 12
$ws.Range("D12").Value = "7.52"
$ws.Range("E12").Value = "  +0.04%  "

# Row 13
$ws.Range("E13").Value = "  -3.68%  "

# Row 14
$ws.Range("D14").Value = "2.905.35"
$ws.Range("E14").Value = "  +0.80%  "

# Row 15
$ws.Range("D15").Value = "2.542.12"
$ws.Range("E15").Value = "  +0.39%  "

# Row 16
$ws.Range("D16").Value = "15.24"
$ws.Range("E16").Value = "  -1.89%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.850"
$ws.Range("E17").Value = "  -0.48%  "

# Row 18
$ws.Range("D18").Value = "42.717.89"
$ws.Range("E18").Value = "  +0.87%  "

# Row 19
$ws.Range("D19").Value = "12.76"
$ws.Range("E19").Value = "  -1.56%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.70"
$ws.Range("E20").Value = "  +4.08%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0958"
$ws.Range("E21").Value = "  -1.55%  "

# Row 22
$ws.Range("D22").Value = "69.52"
$ws.Range("E22").Value = "  -2.61%  "

# Row 23
$ws.Range("D23").Value = "249.59"
$ws.Range("E23").Value = "  -1.43%  "

# Row 24
$ws.Range("E24").Value = "  +1.07%  "

# Row 25
$ws.Range("D25").Value = "2.08"
$ws.Range("E25").Value = "  +2.02%  "

# Row 26
$ws.Range("D26").Value = "26.52"
$ws.Range("E26").Value = "  -1.56%  "

# Row 27
$ws.Range("E27").Value = "  +0.02%  "

# Row 28
$ws.Range("D28").Value = "2.42"
$ws.Range("E28").Value = "  +3.85%  "

# Row 29
$ws.Range("D29").Value = "41.28"
$ws.Range("E29").Value = "  +10.21%  "

# Row 30
$ws.Range("E30").Value = "  +1.34%  "

# Row 31
$ws.Range("D31").Value = "5.94"
$ws.Range("E31").Value = "  +0.20%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "158.00"
$ws.Range("E32").Value = "  +2.13%  "

# Row 33
$ws.Range("B33").Value = "Celestia"
$ws.Range("C33").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D33").Value = "19.31"
$ws.Range("E33").Value = "  +0.30%  "

# Row 34
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "2.14"
$ws.Range("E34").Value = "  +3.17%  "

# Row 35
$ws.Range("D35").Value = "2.68"
$ws.Range("E35").Value = "  +3.26%  "

# Row 36
$ws.Range("D36").Value = "3.29"
$ws.Range("E36").Value = "  +0.58%  "

# Row 37
$ws.Range("D37").Value = "0.0778"
$ws.Range("E37").Value = "  -1.06%  "

# Row 38
$ws.Range("E38").Value = "  -2.40%  "

# Row 39
$ws.Range("E39").Value = "  -0.89%  "

# Row 40
$ws.Range("D40").Value = "23.23"
$ws.Range("E40").Value = "  -3.45%  "

# Row 41
$ws.Range("D41").Value = "2.31"
$ws.Range("E41").Value = "  +14.08%  "

# Row 42
$ws.Range("D42").Value = "0.0304"
$ws.Range("E42").Value = "  +1.13%  "

# Row 43
$ws.Range("E43").Value = "  +0.35%  "

# Row 44
$ws.Range("D44").Value = "3.32"
$ws.Range("E44").Value = "  -2.08%  "

# Row 45
$ws.Range("D45").Value = "3.78"
$ws.Range("E45").Value = "  -2.91%  "

# Row 46
$ws.Range("D46").Value = "2.036.18"
$ws.Range("E46").Value = "  +0.22%  "

# Row 47
$ws.Range("D47").Value = "84.24"
$ws.Range("E47").Value = "  -0.34%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.90"
$ws.Range("E48").Value = "  -1.16%  "

# Row 49
$ws.Range("D49").Value = "74.92"
$ws.Range("E49").Value = "  +2.61%  "

# Row 50
$ws.Range("D50").Value = "105.12"
$ws.Range("E50").Value = "  +3.85%  "

# Row 51
$ws.Range("D51").Value = "2.761.11"
$ws.Range("E51").Value = "  +0.55%  "
